$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every "Price" (D) cell in this sheet holds plain text, even though many
# of the values look like numbers (e.g. "241.80", "0.0960", "0.0000107").
# Assigning such a string straight to .Value lets Excel auto-detect it as
# a number -- silently dropping significant trailing/leading zeros and
# changing the cell type. Mark every edited Price cell as Text ("@") right
# before the write, then restore the default "Normal" style afterwards so
# no stray number-format is left behind (matches the source formatting,
# which never applied an explicit style to these cells).

$ws.Range("D2").Value = "42.013.95"
$ws.Range("D3").Value = "2.216.50"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0960"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "2.549.46"
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "2.203.15"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "41.909.24"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.21%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -6.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.85%  "
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  -9.56%  "
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0305"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D51").Value = "2.423.65"
$ws.Range("E51").Value = "  -1.71%  "
